$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename the header row: "<col>_old" -> "<col>_FV2304", "<col>_new" -> "<col>_FV2310" ---
$headerMap = @{
    1  = "Segmentname_FV2304"
    2  = "Segmentgruppe_FV2304"
    3  = "Segment_FV2304"
    4  = "Datenelement_FV2304"
    5  = "Segment ID_FV2304"
    6  = "Code_FV2304"
    7  = "Qualifier_FV2304"
    8  = "Beschreibung_FV2304"
    9  = "Bedingungsausdruck_FV2304"
    10 = "Bedingung_FV2304"
    11 = "diff"
    12 = "Segmentname_FV2310"
    13 = "Segmentgruppe_FV2310"
    14 = "Segment_FV2310"
    15 = "Datenelement_FV2310"
    16 = "Segment ID_FV2310"
    17 = "Code_FV2310"
    18 = "Qualifier_FV2310"
    19 = "Beschreibung_FV2310"
    20 = "Bedingungsausdruck_FV2310"
    21 = "Bedingung_FV2310"
}

foreach ($col in $headerMap.Keys) {
    $ws.Cells.Item(1, $col).Value = $headerMap[$col]
}

# --- 2. Turn the used range into a real Excel Table (ListObject) so the ---
#        new header names become the table's column names too.
$lo = $ws.ListObjects.Add(1, $ws.Range("A1:U67"), [System.Type]::Missing, 1)
$lo.Name = "Table1"

# --- 3. Freeze the header row (pane split after row 1) ---
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
